$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.996.23"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +0.81%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'3.431.11"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +1.75%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  -0.09%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'408.84"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +1.29%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'128.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -2.87%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.627"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +7.12%  "
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'  -0.12%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.730"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +9.81%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  +22.67%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'42.48"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +1.04%  "
$ws.Range("E11").ClearFormats()
$ws.Range("B12").Value = "'ShibaInu"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C12").ClearFormats()
$ws.Range("D12").Value = "'0.0000217"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +70.36%  "
$ws.Range("E12").ClearFormats()
$ws.Range("B13").Value = "'TRON"
$ws.Range("B13").ClearFormats()
$ws.Range("C13").Value = "'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("C13").ClearFormats()
$ws.Range("D13").Value = "'0.141"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  -0.21%  "
$ws.Range("E13").ClearFormats()
$ws.Range("B14").Value = "'WrappedliquidstakedEther2.0"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("C14").ClearFormats()
$ws.Range("D14").Value = "'3.973.99"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +1.78%  "
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = "'Chainlink"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C15").ClearFormats()
$ws.Range("D15").Value = "'21.33"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +8.28%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'8.89"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +6.21%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'3.434.42"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +1.56%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'12.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +13.70%  "
$ws.Range("E18").ClearFormats()
$ws.Range("E19").Value = "'  +5.77%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'61.880.84"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +0.65%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'394.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +25.44%  "
$ws.Range("E21").ClearFormats()
$ws.Range("E22").Value = "'  +6.36%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'3.20"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +0.63%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'13.24"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  +4.41%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  +3.72%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'32.71"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +11.60%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'8.75"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +6.99%  "
$ws.Range("E27").ClearFormats()
$ws.Range("E28").Value = "'  +0.68%  "
$ws.Range("E28").ClearFormats()
$ws.Range("E29").Value = "'  -0.63%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  +2.10%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  +2.16%  "
$ws.Range("E31").ClearFormats()
$ws.Range("B32").Value = "'InjectiveProtocol"
$ws.Range("B32").ClearFormats()
$ws.Range("C32").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C32").ClearFormats()
$ws.Range("D32").Value = "'43.83"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +5.33%  "
$ws.Range("E32").ClearFormats()
$ws.Range("B33").Value = "'Kaspa"
$ws.Range("B33").ClearFormats()
$ws.Range("C33").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("C33").ClearFormats()
$ws.Range("D33").Value = "'0.171"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +0.01%  "
$ws.Range("E33").ClearFormats()
$ws.Range("B34").Value = "'Cosmos"
$ws.Range("B34").ClearFormats()
$ws.Range("C34").Value = "'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("C34").ClearFormats()
$ws.Range("D34").Value = "'11.90"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +5.31%  "
$ws.Range("E34").ClearFormats()
$ws.Range("E35").Value = "'  +0.00%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.0503"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +5.40%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'53.69"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +3.97%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = "'  -0.06%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  -0.10%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  -0.21%  "
$ws.Range("E40").ClearFormats()
$ws.Range("E41").Value = "'  +6.97%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.313"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +7.57%  "
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'141.82"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +2.38%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  +0.60%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'4.04"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +2.35%  "
$ws.Range("E45").ClearFormats()
$ws.Range("E46").Value = "'  +10.30%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'16.72"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +0.83%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'21.70"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +1.99%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'2.119.95"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +0.25%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  +16.41%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'1.97"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +5.27%  "
$ws.Range("E51").ClearFormats()
